$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 5); the new data set only has 3 rows (2-4)
$ws.Rows.Item(5).Delete()

# Row 2: Resolving-Mac | Cd28 | Cd86 | ECs
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Cd28"
$ws.Range("C2").Value = "Cd86"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.961972333333333
$ws.Range("H2").Value = 14.885917
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.020562
$ws.Range("N2").Value = 0.061686
$ws.Range("O2").Value = 0.0002281281878049052
$ws.Range("P2").Value = 0.0002281281878049052
$ws.Range("Q2").Value = 0.102028075118
$ws.Range("R2").Value = 0.9182526760619999
$ws.Range("S2").Value = 0.0002281281878049052
$ws.Range("T2").Value = 0.0002281281878049052

# Row 3: Resolving-Mac | Cd28 | Cd86 | MuSCs
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Cd28"
$ws.Range("C3").Value = "Cd86"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.961972333333333
$ws.Range("H3").Value = 14.885917
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01154533333333333
$ws.Range("N3").Value = 0.034636
$ws.Range("O3").Value = 0.0001280914293812323
$ws.Range("P3").Value = 0.0001280914293812324
$ws.Range("Q3").Value = 0.0572876245791111
$ws.Range("R3").Value = 0.515588621212
$ws.Range("S3").Value = 0.0001280914293812323
$ws.Range("T3").Value = 0.0001280914293812324

# Row 4: Resolving-Mac | Cd28 | Cd86 | Resolving-Mac
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Cd28"
$ws.Range("C4").Value = "Cd86"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.961972333333333
$ws.Range("H4").Value = 14.885917
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 90.10142766666667
$ws.Range("N4").Value = 270.304283
$ws.Range("O4").Value = 0.9996437803828139
$ws.Range("P4").Value = 0.9996437803828139
$ws.Range("Q4").Value = 447.0807912758345
$ws.Range("R4").Value = 4023.727121482511
$ws.Range("S4").Value = 0.9996437803828139
$ws.Range("T4").Value = 0.9996437803828139
